# "Add table with TCAs" - append a new algorithm entry (FETC) to the
# Algorithms sheet as row 15, and update the sheet view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Algorithms")

# New row of data: Short Name | Long Name | Domain | ... | Reference(s)
$ws.Range("A15").Value = "FETC"
$ws.Range("B15").Value = "Fair and Energyefficient TC"
$ws.Range("C15").Value = "WSN"
$ws.Range("K15").Value = "A Topology Control Protocol for 2D Poisson Distributed Wireless Sensor Networks http://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=5136711 "
$ws.Range("K15").WrapText = $true

# View-state changes captured by the diff: the sheet was re-zoomed and the
# active selection moved to just past the new row.
$excel.ActiveWindow.Zoom = 85
$ws.Range("K16").Select() | Out-Null
